$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# B5, B6, B7 currently all hold "Default ProdBuy".
# Differentiate them into "Default ProdBuy-1", "-2", "-3" respectively.
$ws.Range("B5").Value = "Default ProdBuy-1"
$ws.Range("B6").Value = "Default ProdBuy-2"
$ws.Range("B7").Value = "Default ProdBuy-3"
